$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Update description texts for "CREAR GRUPO" (row 14) and "EDITAR GRUPO" (row 15)
$ws.Range("C14").Value = "En este caso de uso el Director creará un nuevo Grupo que tomará clases en las instalaciones, a dicho Grupo se le asignaran: un nombre, un Maestro, un monto de inscripción y de mensualidad "
$ws.Range("C15").Value = "En este caso de uso el Director cambia la asignación del nombre del grupo, maestro y montos del grupo que haya seleccionado"

# Update effort estimates
$ws.Range("F14").Value = 27
$ws.Range("F15").Value = 24.33
$ws.Range("F30").Value = 10.63
$ws.Range("F32").Value = 4.53

# Update selected cell / view for the sheet (active view)
$ws.Range("F15").Select()
